$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that represents the last
# time the data was refreshed. The automatic update bumped this value
# by one day (45178 -> 45179) for every data row (rows 2 through 533).
$ws.Range("C2:C533").Value = 45179
